$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 321; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "'99999999"
    $cell.Style = "Normal"
}
